# Updated cryptos list on Mon Feb 20 23:41:09 UTC 2023 with GitHub Actions
#
# The "Price" (D) and "Volume(1h)" (E) columns are stored as plain text
# (not numbers) even though many of the Price strings look numeric
# (e.g. "314.96", "1.040"). If we assign those values to Range.Value
# directly, Excel's COM layer auto-detects them as numbers and silently
# normalizes them (dropping trailing zeros, re-parsing "1.040" -> 1.04,
# etc.), which corrupts the text. To keep them as literal text we:
#   1. force the cell to Text format ("@") before assigning,
#   2. assign the literal string,
#   3. clear the formatting again so the cell's style matches the
#      original (unstyled) cells instead of picking up an explicit
#      "Text" number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# --- Row 9 / Row 10 swap (BinanceUSD <-> Polygon) ---
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D9" "1.473"
Set-TextValue "E9" "  -1.26%  "

$ws.Range("B10").Value = "BinanceUSD"
$ws.Range("C10").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D10" "0.9985"
Set-TextValue "E10" "  -0.38%  "

# --- Row 50 / Row 51 swap (Cronos <-> Quant) ---
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "140.02"
Set-TextValue "E50" "  +1.35%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.08374"
Set-TextValue "E51" "  +4.35%  "

# --- Price / Volume refresh for all other rows ---
$updates = [ordered]@{
    'D2'  = '24.777.00'
    'E2'  = '  +1.26%  '
    'D3'  = '1.701.14'
    'E3'  = '  +0.68%  '
    'D4'  = '0.9994'
    'E4'  = '  -0.43%  '
    'D5'  = '314.96'
    'E5'  = '  +0.33%  '
    'E6'  = '  -0.27%  '
    'D7'  = '0.3986'
    'E7'  = '  +2.77%  '
    'D8'  = '0.4036'
    'E8'  = '  +0.40%  '
    'D11' = '53.40'
    'E11' = '  +1.90%  '
    'D12' = '0.08810'
    'E12' = '  +0.75%  '
    'D13' = '26.09'
    'E13' = '  +4.12%  '
    'D14' = '7.547'
    'E14' = '  +0.52%  '
    'D15' = '8.001'
    'E15' = '  +0.07%  '
    'D16' = '0.00001349'
    'E16' = '  +0.15%  '
    'D17' = '1.732.89'
    'E17' = '  +2.80%  '
    'D18' = '95.83'
    'E18' = '  -2.48%  '
    'D19' = '0.07188'
    'E19' = '  +1.45%  '
    'E20' = '  +4.37%  '
    'D21' = '7.347'
    'E21' = '  +1.41%  '
    'E22' = '  -0.28%  '
    'E23' = '  +1.01%  '
    'D24' = '24.779.55'
    'E24' = '  +1.26%  '
    'D25' = '2.374'
    'E25' = '  +0.83%  '
    'D26' = '2.929'
    'E26' = '  -1.17%  '
    'D27' = '23.18'
    'E27' = '  +2.22%  '
    'D28' = '6.164'
    'E28' = '  +18.22%  '
    'D29' = '161.73'
    'E29' = '  -0.69%  '
    'D30' = '8.635'
    'E30' = '  -1.16%  '
    'D31' = '144.22'
    'E31' = '  +5.39%  '
    'D32' = '2.442'
    'E32' = '  +24.88%  '
    'D33' = '1.909.73'
    'E33' = '  +2.08%  '
    'D34' = '0.08646'
    'E34' = '  -1.99%  '
    'D35' = '7.342'
    'E35' = '  -1.07%  '
    'D36' = '0.03179'
    'E36' = '  +9.66%  '
    'D37' = '1.040'
    'E37' = '  +0.67%  '
    'D38' = '0.2840'
    'E38' = '  +0.73%  '
    'D39' = '10.79'
    'E39' = '  +0.21%  '
    'D40' = '0.09447'
    'E40' = '  +3.61%  '
    'D41' = '0.8301'
    'E41' = '  +4.83%  '
    'D42' = '14.23'
    'E42' = '  +0.21%  '
    'D43' = '1.479'
    'E43' = '  +1.66%  '
    'D44' = '17.63'
    'E44' = '  +5.49%  '
    'D45' = '2.701'
    'E45' = '  +4.01%  '
    'D46' = '0.7438'
    'E46' = '  +2.85%  '
    'D47' = '4.210'
    'E47' = '  +0.21%  '
    'D48' = '1.398'
    'E48' = '  +4.18%  '
    'E49' = '  -0.23%  '
}

foreach ($ref in $updates.Keys) {
    Set-TextValue $ref $updates[$ref]
}
